$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '48.014.54'
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').Value = '2.510.39'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.63'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.61'
$ws.Range('E6').Value = '  +4.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.524'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.548'
$ws.Range('E9').Value = '  +0.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.47'
$ws.Range('E10').Value = '  +5.91%  '
$ws.Range('E11').Value = '  +0.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.125'
$ws.Range('E12').Value = '  +0.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.70'
$ws.Range('E13').Value = '  +2.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.27'
$ws.Range('E14').Value = '  +1.18%  '
$ws.Range('D15').Value = '2.899.39'
$ws.Range('E15').Value = '  +0.74%  '
$ws.Range('D16').Value = '2.512.09'
$ws.Range('E16').Value = '  +0.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.853'
$ws.Range('E17').Value = '  +0.58%  '
$ws.Range('D18').Value = '47.851.15'
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.32'
$ws.Range('E19').Value = '  +4.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.67'
$ws.Range('E20').Value = '  +1.36%  '
$ws.Range('E21').Value = '  +0.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.76'
$ws.Range('E22').Value = '  +14.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.86'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '248.14'
$ws.Range('E24').Value = '  -1.36%  '
$ws.Range('E25').Value = '  -0.85%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.96'
$ws.Range('E27').Value = '  -0.99%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  +2.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.10'
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.87'
$ws.Range('E32').Value = '  +0.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.21'
$ws.Range('E33').Value = '  +2.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.37'
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0790'
$ws.Range('E35').Value = '  +0.66%  '
$ws.Range('E36').Value = '  +0.39%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.98'
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.72'
$ws.Range('E38').Value = '  +1.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.98'
$ws.Range('E39').Value = '  -0.46%  '
$ws.Range('E40').Value = '  +0.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.59'
$ws.Range('E41').Value = '  +6.90%  '
$ws.Range('E42').Value = '  -1.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '119.38'
$ws.Range('E43').Value = '  -2.06%  '
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('D45').Value = '2.001.44'
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.07'
$ws.Range('E46').Value = '  +3.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.04'
$ws.Range('E47').Value = '  -3.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.88'
$ws.Range('E48').Value = '  +3.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.08'
$ws.Range('E49').Value = '  -1.31%  '
$ws.Range('E50').Value = '  -0.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '56.79'
$ws.Range('E51').Value = '  +3.62%  '
